$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update single values in N15 and N16 ---
$ws.Range("N15").Value2 = 129286.935
$ws.Range("N16").Value2 = 35684.474000000002

# --- Swap the content (label + data) of rows 17 and 18, keeping each ---
# --- row's own formatting (styles) in place.                        ---

# Capture current ("before") content of row 17 and row 18.
$label17 = $ws.Range("B17").Value2
$label18 = $ws.Range("B18").Value2
$data17 = $ws.Range("C17:N17").Value2
$data18 = $ws.Range("C18:N18").Value2

# Row 17 becomes what row 18 used to hold.
$ws.Range("B17").Value2 = $label18
$ws.Range("C17:N17").Value2 = $data18

# Row 18 becomes what row 17 used to hold (i.e. empty data cells).
$ws.Range("B18").Value2 = $label17
$ws.Range("C18:N18").Value2 = $data17

# --- Row 21 formatting: number format now matches the one used on ---
# --- row 18 (style index 7) instead of the bespoke/duplicate style ---
# --- that gets dropped from the workbook's style table entirely.  ---
$ws.Range("C21:N21").NumberFormat = $ws.Range("C18").NumberFormat
